$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: MCH228-1 ---
$ws.Range("A2").Value = "MCH228-1"
$ws.Range("C2").Value = "CONSTITUTION, MINUTES, REPORT, CORRESPONDENCE, CIRCULARS, FINANCE, CONFERENCE, EDUCATION"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24F | GRAP COUNT NUMER: NONE"
$ws.Range("H2").Value = ""

# --- Row 3: MCH228-2 ---
$ws.Range("A3").Value = "MCH228-2"
$ws.Range("C3").Value = "EDUCATION, MEMORANDUM, PRESS STATEMENTS, ADDRESSES, MESSAGES, LITERACY, PAPER CUTTINGS"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 24F | GRAP COUNT NUMER: NONE"
$ws.Range("H3").Value = ""

# --- Formatting to match the rest of the data rows: 10pt Calibri, theme text color ---
$dataRange = $ws.Range("A2:H3")
$dataRange.Font.Name = "Calibri"
$dataRange.Font.ThemeColor = 1

# Columns E/F/G need the same font/color treatment applied individually as well
# (already covered by the A2:H3 range above)

# F column (extentAndMedium) additionally carries the alignment-applied flag
$ws.Range("F2:F3").WrapText = $false

# --- Row heights to match authored rows ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

# --- View: freeze header row, select the new data rows ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:J3").Select()
